# Update odds values in Sheet1 to match the 2024-10-31 FlashScore refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (AS Roma - Torino)
$ws.Range("G2").Value = 1.8
$ws.Range("H2").Value = 3.6
$ws.Range("I2").Value = 4.33
$ws.Range("J2").Value = 2.5
$ws.Range("N2").Value = 9.5
$ws.Range("O2").Value = 1.36
$ws.Range("P2").Value = 3.2
$ws.Range("Z2").Value = 15
$ws.Range("AW2").Value = 6

# Row 3 (Como - Lazio)
$ws.Range("Q3").Value = 2.04
$ws.Range("R3").Value = 1.86

# Row 5 (Sarmiento Junin - Independiente)
$ws.Range("M5").Value = 1.14
$ws.Range("N5").Value = 5.5

# Row 7 (Grasshoppers - Lugano)
$ws.Range("G7").Value = 3.3
$ws.Range("H7").Value = 3.4
$ws.Range("I7").Value = 2.1
$ws.Range("J7").Value = 3.75
$ws.Range("L7").Value = 2.75
$ws.Range("X7").Value = 17
$ws.Range("Y7").Value = 12
$ws.Range("Z7").Value = 34
$ws.Range("AB7").Value = 29
$ws.Range("AI7").Value = 11
$ws.Range("AL7").Value = 17
$ws.Range("AO7").Value = 17
$ws.Range("AP7").Value = 23
$ws.Range("AR7").Value = 67
$ws.Range("AY7").Value = 21
$ws.Range("AZ7").Value = 41

# Row 8 (Servette - Luzern)
$ws.Range("H8").Value = 3.7
$ws.Range("N8").Value = 15
$ws.Range("O8").Value = 1.18
$ws.Range("P8").Value = 4.5
$ws.Range("Q8").Value = 1.62
$ws.Range("R8").Value = 2.25
$ws.Range("S8").Value = 1.3
$ws.Range("T8").Value = 3.4
$ws.Range("U8").Value = 1.57
$ws.Range("V8").Value = 2.25
$ws.Range("W8").Value = 9.5
$ws.Range("X8").Value = 10
$ws.Range("Z8").Value = 15
$ws.Range("AB8").Value = 21
$ws.Range("AC8").Value = 15
$ws.Range("AF8").Value = 41
$ws.Range("AG8").Value = 126
$ws.Range("AH8").Value = 15
$ws.Range("AJ8").Value = 13
$ws.Range("AN8").Value = 4
$ws.Range("AS8").Value = 101
$ws.Range("AT8").Value = 3.4
$ws.Range("AU8").Value = 7.5
$ws.Range("AX8").Value = 21
$ws.Range("AZ8").Value = 67
